$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "[Deep Learning with Python] 2판 번역 완료!"
$ws.Range("E12").Value = "https://tensorflow.blog/2022/06/11/deep-learning-with-python-2%ed%8c%90-%eb%b2%88%ec%97%ad-%ec%99%84%eb%a3%8c/"

$ws.Range("D24").Value = "자랑스럽고도 부끄러운"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222768616651"

$ws.Range("D51").Value = "[윈도우11] 어떤 프로그램의 바로 가기를 최소화된 창으로 실행되게 하려면?"
$ws.Range("E51").Value = "https://bskyvision.com/1292"
